$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "TC001"
$ws.Range("B8").Value = "Creating the TEAM Workspace"
$ws.Range("C8").Value = "TEAM Workspace Should be created successfully and approved"
$ws.Range("D8").Value = "TEAM workspace is created successfully and approved"
$ws.Range("E8").Value = "Pass"
